$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp update
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 13:30"

# Row 7
$ws.Range("B7").Value = 1431635
$ws.Range("C7").Value = 16319
$ws.Range("D7").Value = 1085608
$ws.Range("E7").Value = 321392
$ws.Range("G7").Value = 269
$ws.Range("H7").Value = 24635

# Row 16
$ws.Range("B16").Value = 539670
$ws.Range("C16").Value = 5039
$ws.Range("D16").Value = 434676
$ws.Range("E16").Value = 73960
$ws.Range("G16").Value = 322
$ws.Range("H16").Value = 31034

# Row 21
$ws.Range("B21").Value = 374734
$ws.Range("C21").Value = 1003
$ws.Range("D21").Value = 298300
$ws.Range("E21").Value = 66524
$ws.Range("G21").Value = 11
$ws.Range("H21").Value = 9910

# Row 26
$ws.Range("B26").Value = 324077
$ws.Range("C26").Value = 618
$ws.Range("D26").Value = 308020
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 6673

# Row 33
$ws.Range("B33").Value = 186254
$ws.Range("C33").Value = 3400
$ws.Range("D33").Value = 134395
$ws.Range("E33").Value = 45863
$ws.Range("G33").Value = 65
$ws.Range("H33").Value = 5996

# Row 38
$ws.Range("B38").Value = 139129
$ws.Range("C38").Value = 3093
$ws.Range("D38").Value = 96609
$ws.Range("E38").Value = 41755
$ws.Range("G38").Value = 8
$ws.Range("H38").Value = 765

# Row 54
$ws.Range("B54").Value = 88909
$ws.Range("C54").Value = 619
$ws.Range("D54").Value = 80503
$ws.Range("E54").Value = 7469
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = 937

# Row 56
$ws.Range("A56").Value = "Suiza"
$ws.Range("B56").Value = 86167
$ws.Range("C56").Value = 3008
$ws.Range("D56").Value = 54600
$ws.Range("E56").Value = 29423
$ws.Range("G56").Value = 6
$ws.Range("H56").Value = 2144

# Row 57
$ws.Range("A57").Value = "China"
$ws.Range("B57").Value = 85704
$ws.Range("C57").Value = 19
$ws.Range("D57").Value = 80812
$ws.Range("E57").Value = 258
$ws.Range("H57").Value = 4634

# Row 68
$ws.Range("B68").Value = 52910
$ws.Range("C68").Value = 384
$ws.Range("D68").Value = 46139
$ws.Range("E68").Value = 5658
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 1113

# Row 70
$ws.Range("B70").Value = 50906
$ws.Range("C70").Value = 957
$ws.Range("D70").Value = 27832
$ws.Range("E70").Value = 22328
$ws.Range("G70").Value = 14
$ws.Range("H70").Value = 746

# Row 77
$ws.Range("B77").Value = 40357
$ws.Range("C77").Value = 70
$ws.Range("D77").Value = 33790
$ws.Range("E77").Value = 5068
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 1499

# Row 100
$ws.Range("B100").Value = 15459
$ws.Range("C100").Value = 27
$ws.Range("D100").Value = 13922
$ws.Range("E100").Value = 1217
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 320

# Row 101
$ws.Range("A101").Value = "Eslovenia"
$ws.Range("B101").Value = 14473
$ws.Range("C101").Value = 794
$ws.Range("D101").Value = 6572
$ws.Range("E101").Value = 7709
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 192

# Row 102
$ws.Range("A102").Value = "Finlandia"
$ws.Range("B102").Value = 13849
$ws.Range("C102").Value = 294
$ws.Range("D102").Value = 9100
$ws.Range("E102").Value = 4398
$ws.Range("H102").Value = 351

# Row 103
$ws.Range("A103").Value = "Sudan"
$ws.Range("B103").Value = 13724
$ws.Range("D103").Value = 6764
$ws.Range("E103").Value = 6124
$ws.Range("H103").Value = 836

# Row 126
$ws.Range("B126").Value = 5685
$ws.Range("C126").Value = 60
$ws.Range("D126").Value = 3457
$ws.Range("E126").Value = 2215

# Row 131
$ws.Range("B131").Value = 5262
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 4996
$ws.Range("E131").Value = 161

# Row 138
$ws.Range("A138").Value = "Malta"
$ws.Range("B138").Value = 4871
$ws.Range("C138").Value = 134
$ws.Range("D138").Value = 3282
$ws.Range("E138").Value = 1543
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 46

# Row 139
$ws.Range("A139").Value = "Republica de Africa Central"
$ws.Range("B139").Value = 4856
$ws.Range("D139").Value = 1924
$ws.Range("E139").Value = 2870
$ws.Range("H139").Value = 62
